# "filtering non dispersed data"
#
# The workbook originally has a single sheet ("Sheet1") containing a table of
# book-valuation statistics grouped by "Especialidad"/"Libro". The last two
# groups (rows 20-31: "Eje Transversal" / Libro 68 and Libro 69) are considered
# "dispersed" data and are filtered out into their own situation:
#   * The original sheet becomes an empty sheet renamed "Sheet".
#   * A new sheet named "tabla de datos dispersos" is created holding only the
#     first three groups (rows 1-19, i.e. Economia/Libro 42, Libro 18, Libro 25).
#
# We duplicate the original sheet first (so the duplicate - not a brand new
# blank sheet - keeps every row/style/merge byte-for-byte) and trim the
# duplicate down to rows 1-19. The original sheet object is then wiped to
# become the new blank "Sheet" tab. This keeps the low sheetId (1) on the
# first tab ("Sheet") and the higher sheetId (2) on the new data tab, matching
# how Excel numbers sheets in creation order.

$wb = $excel.ActiveWorkbook
$wsOrig = $wb.Worksheets.Item(1)

# Duplicate the whole sheet right after itself - a full-sheet copy carries
# over styles/merges exactly (no border-splitting side effects like a
# Range.Copy into a fresh sheet would cause).
$wsOrig.Copy($null, $wsOrig)

$wsBlank = $wb.Worksheets.Item(1)
$wsData = $wb.Worksheets.Item(2)

# Trim the copy down to the non-dispersed rows only (header + first 3 groups).
$wsData.Range("A20:L31").EntireRow.Delete()
$wsData.Name = "tabla de datos dispersos"

# Wipe the original sheet clean (drop the merges first so they don't linger
# once the cells are cleared) and rename it.
$wsBlank.Cells.UnMerge()
$wsBlank.Cells.Clear()
$wsBlank.Name = "Sheet"

# Keep the first tab ("Sheet") selected/active, same as before the edit.
$wsBlank.Activate()
